$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Significant"

for ($row = 2; $row -le 10; $row++) {
    $pValue = $ws.Cells.Item($row, 2).Value2
    $ws.Cells.Item($row, 4).Value = ($pValue -lt 0.05)
}
